$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text/string storage for numeric-looking price cells so Excel
# does not auto-convert them to numbers (which would drop formatting
# like trailing zeros, e.g. "35.20" -> 35.2).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.901.84"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.329.22"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "605.88"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "143.22"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.327.76"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +1.80%  "
$ws.Range("D11").Value = "5.58"
$ws.Range("E11").Value = "  +4.40%  "
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("E13").Value = "  +0.52%  "
$ws.Range("D14").Value = "35.20"
$ws.Range("E14").Value = "  +2.19%  "
$ws.Range("D15").Value = "3.884.50"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "3.335.85"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "64.009.39"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "6.88"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "482.99"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "14.12"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "14.05"
$ws.Range("E24").Value = "  +6.32%  "
$ws.Range("D25").Value = "85.04"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("D28").Value = "8.31"
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("E31").Value = "  +2.24%  "
$ws.Range("D32").Value = "28.94"
$ws.Range("E32").Value = "  +4.83%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "6.10"
$ws.Range("E36").Value = "  +3.09%  "
$ws.Range("D37").Value = "0.0₃0754"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("D38").Value = "52.48"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("D39").Value = "0.0400"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "435.64"
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "3.136.39"
$ws.Range("E41").Value = "  +5.15%  "
$ws.Range("E42").Value = "  +7.18%  "
$ws.Range("D43").Value = "2.77"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "8.38"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").Value = "36.98"
$ws.Range("E47").Value = "  +8.95%  "
$ws.Range("D48").Value = "26.53"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "124.78"
$ws.Range("E51").Value = "  +3.01%  "
